$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Range("A2").Value = "Drew Koecher"
$ws.Activate()
$ws.Range("E16").Select()
